$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Department columns (P:T) for rows 2-4 to add the prior year department
foreach ($r in 2..4) {
    $ws.Cells.Item($r, 16).Value = 3      # P: Department ID

    $qCell = $ws.Cells.Item($r, 17)       # Q: Department Name
    $qCell.NumberFormat = "@"             # force text storage so "2" stays a string
    $qCell.Value = "2"
    $qCell.Style = "Normal"               # restore default style (no number format override)

    $ws.Cells.Item($r, 18).Value = 100    # R: Default Total Hours
    $ws.Cells.Item($r, 19).Value = 21     # S: Core Skills Percentage
    $ws.Cells.Item($r, 20).Value = 79     # T: Soft Skills Percentage
}
